# Append new order-line rows (30-39) to the active worksheet, matching
# the VendorManager / CLI dynamic bot calling update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Smith Quality Eggs LLC", "Eggs (In-Shell)", "1", "$0.00", "$0.00"),
    @("", "Flour - Millers Choice", "1", "$0.00", "$0.00"),
    @("", "Goat Cheese", "7", "$0.00", "$0.00"),
    @("", "Mustard - Honey", "3", "$0.00", "$0.00"),
    @("", "Nuts - Almonds (sliced)", "1", "$106.46", "$106.46"),
    @("Casa", "Nuts - Pine", "3", "$117.37", "$352.11"),
    @("", "Quinoa - Tri Color", "3", "$58.31", "$174.93"),
    @("", "Sugar - Extra Fine", "1", "$0.00", "$0.00"),
    @("", "Nuts-Almonds Blanched (Slivered)", "1", "$0.00", "$0.00"),
    @("PERF", "Vegan Egg", "1", "$99.59", "$99.59")
)

$startRow = 30
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]

    # The source sheet stores every value (including quantities and
    # currency amounts) as plain text, so force text formatting on the
    # whole row before writing so Excel does not auto-coerce these
    # numeric-looking strings into real numbers.
    $rowRange = $ws.Range("A" + $r + ":E" + $r)
    $rowRange.NumberFormat = "@"

    $ws.Range("A$r").Value = $rowValues[0]
    $ws.Range("B$r").Value = $rowValues[1]
    $ws.Range("C$r").Value = $rowValues[2]
    $ws.Range("D$r").Value = $rowValues[3]
    $ws.Range("E$r").Value = $rowValues[4]
}
